# ContactTemplate.xlsx update: "added parse to excel method"
# Rebuilds sheet1 with a new header layout (Name/Email/Phone number/Contact Type/
# Is Company/Title/Customer type/Zip/Street/City/Country/VAT number) plus a
# sample/test data row, a hyperlink on the Email cell, and trimmed column set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate: drop every existing column (and their widths,
# the old headers, and the now-unused shared strings "Car number"/"Company")
# so the sheet can be rebuilt from scratch exactly like the target layout.
$ws.Range("A1:O1").EntireColumn.Delete()

# --- Column widths (characters) for the columns that keep a custom width ---
$ws.Columns.Item(1).ColumnWidth  = 11.592447916666666   # A -> 12.42578125
$ws.Columns.Item(2).ColumnWidth  = 16.022135416666668   # B -> 16.85546875
$ws.Columns.Item(3).ColumnWidth  = 15.166666666666666   # C -> 16
$ws.Columns.Item(4).ColumnWidth  = 13.451822916666666   # D -> 14.28515625
$ws.Columns.Item(5).ColumnWidth  = 10.451822916666666   # E -> 11.28515625
$ws.Columns.Item(6).ColumnWidth  = 13.736979166666666   # F -> 14.5703125
$ws.Columns.Item(7).ColumnWidth  = 15.877604166666666   # G -> 16.7109375
$ws.Columns.Item(12).ColumnWidth = 12.736979166666666   # L -> 13.5703125

# --- Header row (row 1) ---
# Values are entered in this particular order so the shared-string table is
# built up in the same sequence as the target workbook.
$ws.Range("A1").Value = "Name"
$ws.Range("C1").Value = "Phone number"
$ws.Range("D1").Value = "Contact Type"
$ws.Range("F1").Value = "Title (mr,mrs)"
$ws.Range("H1").Value = "Zip"
$ws.Range("I1").Value = "Street"
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "Country"
$ws.Range("G1").Value = "Customer type"
$ws.Range("B1").Value = "Email"
$ws.Range("L1").Value = "VAT number"

# --- Sample/test data row (row 2) ---
$ws.Range("A2").Value = "test name"
$ws.Range("C2").Value = 3809123123
$ws.Range("K2").Value = "ewe"
$ws.Range("F2").Value = "mr"
$ws.Range("D2").Value = "customer"
$ws.Range("G2").Value = "qwerty"
$ws.Range("H2").Value = 31000
$ws.Range("I2").Value = 232
$ws.Range("J2").Value = "qwe"
$ws.Range("L2").Value = 123456789

# Email cell carries a mailto hyperlink styled with the built-in Hyperlink style
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:wew@mail.com", "", "", "wew@mail.com")

$ws.Range("E1").Value = "Is Company"
# Force "true" to be stored as text (not boolean) by leading with an
# apostrophe, same trick used in the Excel UI.
$ws.Range("E2").Value = "'true"

# --- Selection shown when the sheet is reopened ---
$ws.Range("M6").Select()
